# Logged Week 16 and performed season sim from Week 17
# Update the "R" row (row 3) totals on both the OFF and DEF sheets.

$wb = $excel.ActiveWorkbook

# --- OFF sheet ---
$wsOff = $wb.Worksheets.Item("OFF")
$wsOff.Range("B3").Value = 227
$wsOff.Range("C3").Value = 159
$wsOff.Range("D3").Value = 56
$wsOff.Range("E3").Value = 26
$wsOff.Range("F3").Value = 2
$wsOff.Range("G3").Value = 7

# --- DEF sheet ---
$wsDef = $wb.Worksheets.Item("DEF")
$wsDef.Range("B3").Value = 529
$wsDef.Range("C3").Value = 383
$wsDef.Range("D3").Value = 127
$wsDef.Range("E3").Value = 62
